# Fill in the previously-blank "Preconditions", "Method Inputs" and
# "Expected Result" columns for the first two unit-test rows (test cases
# for the __init__ method) of the demonstration test plan.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Test Case 1 (__init__ / "Attribute set to input value.")
$ws.Range("E7").Value = "none"
$ws.Range("F7").Value = 'name="ISD" department=COMPUTER_SCIENCE credit_hours=6'
$ws.Range("G7").Value = "Client object created with expected attribute values"

# Test Case 2 (__init__ / "Exception raised when name is blank")
$ws.Range("E8").Value = "none"
$ws.Range("F8").Value = 'name="ISD" department=COMPUTER_SCIENCE credit_hours=6'
$ws.Range("G8").Value = "ValueError"
